# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh values to the FFXIV crafting-log sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 255
$ws.Range("I4").Value = 255
$ws.Range("K4").Value = 255
$ws.Range("M4").Value = -141
$ws.Range("H40").Value = 8093.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8093.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8093.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8443.5
$ws.Range("H88").Value = 2062.6667
$ws.Range("I88").Value = 1871
$ws.Range("K88").Value = 1871
$ws.Range("M88").Value = -1465
$ws.Range("H91").Value = 2062.6667
$ws.Range("I91").Value = 1871
$ws.Range("K91").Value = 1871
$ws.Range("M91").Value = -467
$ws.Range("H113").Value = 7621.3335
$ws.Range("I113").Value = 7448.5
$ws.Range("K113").Value = 7448.5
$ws.Range("M113").Value = -4194.5
$ws.Range("H125").Value = 879.75
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H132").Value = 38181.82
$ws.Range("I132").Value = 55310.684
$ws.Range("K132").Value = 165932.052
$ws.Range("M132").Value = -163402.052
$ws.Range("H137").Value = 1394.2222
$ws.Range("I137").Value = 1362.25
$ws.Range("J137").Value = 1650
$ws.Range("K137").Value = 4086.75
$ws.Range("L137").Value = 4950
$ws.Range("M137").Value = -1536.75
$ws.Range("N137").Value = -10050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 15292.9
$ws.Range("I2").Value = 21116.5
$ws.Range("J2").Value = 6557.5
$ws.Range("K2").Value = 21116.5
$ws.Range("L2").Value = 6557.5
$ws.Range("M2").Value = -21003.5
$ws.Range("N2").Value = -6783.5
$ws.Range("H32").Value = 1765.9722
$ws.Range("I32").Value = 1865.6765
$ws.Range("J32").Value = 71
$ws.Range("K32").Value = 1865.6765
$ws.Range("L32").Value = 71
$ws.Range("M32").Value = -1578.6765
$ws.Range("N32").Value = -645
$ws.Range("H116").Value = 15292.9
$ws.Range("I116").Value = 21116.5
$ws.Range("J116").Value = 6557.5
$ws.Range("K116").Value = 21116.5
$ws.Range("L116").Value = 6557.5
$ws.Range("M116").Value = -18822.5
$ws.Range("N116").Value = -11145.5
$ws.Range("H122").Value = 1981.8889
$ws.Range("I122").Value = 1729.625
$ws.Range("K122").Value = 5188.875
$ws.Range("M122").Value = -2738.875
$ws.Range("H132").Value = 12199251
$ws.Range("I132").Value = 2733.697
$ws.Range("J132").Value = 62509884
$ws.Range("K132").Value = 8201.091
$ws.Range("L132").Value = 187529652
$ws.Range("M132").Value = -5671.091
$ws.Range("N132").Value = -187534712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 15292.9
$ws.Range("I3").Value = 21116.5
$ws.Range("J3").Value = 6557.5
$ws.Range("K3").Value = 21116.5
$ws.Range("L3").Value = 6557.5
$ws.Range("M3").Value = -21002.5
$ws.Range("N3").Value = -6785.5
$ws.Range("H22").Value = 568.75
$ws.Range("I22").Value = 542.5
$ws.Range("J22").Value = 574
$ws.Range("K22").Value = 542.5
$ws.Range("L22").Value = 574
$ws.Range("M22").Value = -369.5
$ws.Range("N22").Value = -920
$ws.Range("H107").Value = 31083.766
$ws.Range("I107").Value = 1667.875
$ws.Range("K107").Value = 1667.875
$ws.Range("M107").Value = 252.125
$ws.Range("H134").Value = 19233730
$ws.Range("I134").Value = 20836208
$ws.Range("K134").Value = 62508624
$ws.Range("M134").Value = -62506089

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2200.44
$ws.Range("I31").Value = 1262.9333
$ws.Range("K31").Value = 1262.9333
$ws.Range("M31").Value = -967.9332999999999
$ws.Range("H34").Value = 2200.44
$ws.Range("I34").Value = 1262.9333
$ws.Range("K34").Value = 1262.9333
$ws.Range("M34").Value = -1060.9333
$ws.Range("H94").Value = 595.93335
$ws.Range("I94").Value = 1789
$ws.Range("J94").Value = 510.7143
$ws.Range("K94").Value = 1789
$ws.Range("L94").Value = 510.7143
$ws.Range("M94").Value = -1338
$ws.Range("N94").Value = -1412.7143
$ws.Range("H105").Value = 1499.5
$ws.Range("I105").Value = 1499
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 1499
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 248
$ws.Range("N105").Value = -4994
$ws.Range("H122").Value = 29977.334
$ws.Range("J122").Value = 84733.336
$ws.Range("L122").Value = 254200.008
$ws.Range("N122").Value = -259100.008
$ws.Range("H132").Value = 1685.9333
$ws.Range("I132").Value = 1724.6666
$ws.Range("J132").Value = 1531
$ws.Range("K132").Value = 5173.9998
$ws.Range("L132").Value = 4593
$ws.Range("M132").Value = -2643.9998
$ws.Range("N132").Value = -9653
$ws.Range("H134").Value = 2305.5833
$ws.Range("I134").Value = 2241
$ws.Range("K134").Value = 6723
$ws.Range("M134").Value = -4188

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1587.5
$ws.Range("I9").Value = 1587.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 4762.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -4538.5
$ws.Range("N9").ClearContents()
$ws.Range("H25").Value = 5270.7144
$ws.Range("I25").Value = 2548.75
$ws.Range("K25").Value = 7646.25
$ws.Range("M25").Value = -7477.25
$ws.Range("H30").Value = 5270.7144
$ws.Range("I30").Value = 2548.75
$ws.Range("K30").Value = 7646.25
$ws.Range("M30").Value = -7544.25
$ws.Range("H49").Value = 717
$ws.Range("I49").Value = 717
$ws.Range("K49").Value = 2151
$ws.Range("M49").Value = -1995
$ws.Range("H69").Value = 2383.1667
$ws.Range("J69").Value = 2474.75
$ws.Range("L69").Value = 7424.25
$ws.Range("N69").Value = -9046.25
$ws.Range("H72").Value = 2383.1667
$ws.Range("J72").Value = 2474.75
$ws.Range("L72").Value = 22272.75
$ws.Range("N72").Value = -30384.75
$ws.Range("H81").Value = 7648.25
$ws.Range("I81").Value = 7198.1665
$ws.Range("J81").Value = 8998.5
$ws.Range("K81").Value = 21594.4995
$ws.Range("L81").Value = 26995.5
$ws.Range("M81").Value = -20471.4995
$ws.Range("N81").Value = -29241.5
$ws.Range("H84").Value = 7648.25
$ws.Range("I84").Value = 7198.1665
$ws.Range("J84").Value = 8998.5
$ws.Range("K84").Value = 64783.4985
$ws.Range("L84").Value = 80986.5
$ws.Range("M84").Value = -59167.4985
$ws.Range("N84").Value = -92218.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2694.5
$ws.Range("I80").Value = 2694.5
$ws.Range("K80").Value = 2694.5
$ws.Range("M80").Value = -1696.5
$ws.Range("H83").Value = 2694.5
$ws.Range("I83").Value = 2694.5
$ws.Range("K83").Value = 13472.5
$ws.Range("M83").Value = -8480.5
$ws.Range("H122").Value = 4280
$ws.Range("I122").Value = 4413.8184
$ws.Range("J122").Value = 2808
$ws.Range("K122").Value = 13241.4552
$ws.Range("L122").Value = 8424
$ws.Range("M122").Value = -10791.4552
$ws.Range("N122").Value = -13324
$ws.Range("H126").Value = 2826.9092
$ws.Range("I126").Value = 2899
$ws.Range("J126").Value = 2785.7144
$ws.Range("K126").Value = 8697
$ws.Range("L126").Value = 8357.143199999999
$ws.Range("M126").Value = -6227
$ws.Range("N126").Value = -13297.1432
$ws.Range("H132").Value = 2697.6191
$ws.Range("I132").Value = 2268.6924
$ws.Range("J132").Value = 3394.625
$ws.Range("K132").Value = 6806.0772
$ws.Range("L132").Value = 10183.875
$ws.Range("M132").Value = -4276.0772
$ws.Range("N132").Value = -15243.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3332.4
$ws.Range("I22").Value = 1465
$ws.Range("K22").Value = 1465
$ws.Range("M22").Value = -1170
$ws.Range("H27").Value = 3332.4
$ws.Range("I27").Value = 1465
$ws.Range("K27").Value = 1465
$ws.Range("M27").Value = -1358
$ws.Range("H40").Value = 3497.7856
$ws.Range("I40").Value = 2624.4119
$ws.Range("K40").Value = 2624.4119
$ws.Range("M40").Value = -2488.4119
$ws.Range("H43").Value = 7201.1
$ws.Range("J43").Value = 7201.1
$ws.Range("L43").Value = 7201.1
$ws.Range("N43").Value = -7587.1
$ws.Range("H132").Value = 2650
$ws.Range("I132").Value = 2388.889
$ws.Range("K132").Value = 7166.667
$ws.Range("M132").Value = -4636.667
$ws.Range("H136").Value = 35718532
$ws.Range("I136").Value = 3125.1
$ws.Range("K136").Value = 9375.299999999999
$ws.Range("M136").Value = -6825.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 5987.875
$ws.Range("I113").Value = 1999.3334
$ws.Range("J113").Value = 8381
$ws.Range("K113").Value = 5998.0002
$ws.Range("L113").Value = 25143
$ws.Range("M113").Value = -3828.0002
$ws.Range("N113").Value = -29483
$ws.Range("H122").Value = 2382.5293
$ws.Range("I122").Value = 2513.6667
$ws.Range("J122").Value = 2067.8
$ws.Range("K122").Value = 7541.000100000001
$ws.Range("L122").Value = 6203.400000000001
$ws.Range("M122").Value = -5091.000100000001
$ws.Range("N122").Value = -11103.4
$ws.Range("H132").Value = 3623.625
$ws.Range("I132").Value = 1798.4
$ws.Range("J132").Value = 6665.6665
$ws.Range("K132").Value = 5395.200000000001
$ws.Range("L132").Value = 19996.9995
$ws.Range("M132").Value = -2865.200000000001
$ws.Range("N132").Value = -25056.9995
